$d = $word.ActiveDocument

# Replace the merge-field property name "venue_name" with "external_short_name"
# in "<<caseManagementLocation.venue_name>>..." so the template now reads
# "<<caseManagementLocation.external_short_name>>..."
$rng = $d.Content
$found = $rng.Find.Execute(
    "caseManagementLocation",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0
)
if ($found) {
    $rng.ParagraphFormat.LineSpacingRule = 0
}

$d.Content.Find.Execute(
    ".venue_name",
    $true, $false, $false, $false, $false,
    $true, 1, $false, ".external_short_name", 2
)
